# The workbook is already open as $excel.ActiveWorkbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Permitir scroll nas activitys" task (row 9) as done, matching
# the "FEITO" marker already used in C2, C3 and C8.
$ws.Range("C9").Value = "FEITO"

# Move/save the active selection to A6 (matches the saved cursor position
# recorded in the worksheet's sheetView).
$ws.Range("A6").Select()
